$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 27778296
$ws.Range("J46").Value = 700
$ws.Range("L46").Value = 2100
$ws.Range("N46").Value = -2338
$ws.Range("H58").Value = 2030.3077
$ws.Range("J58").Value = 2909.4119
$ws.Range("L58").Value = 8728.235700000001
$ws.Range("N58").Value = -9028.235700000001
$ws.Range("H59").Value = 1500
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1500
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 4500
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -5614
$ws.Range("H60").Value = 27778296
$ws.Range("J60").Value = 700
$ws.Range("L60").Value = 2100
$ws.Range("N60").Value = -3068
$ws.Range("H70").Value = 1133.919
$ws.Range("I70").Value = 1300.8
$ws.Range("J70").Value = 786.25
$ws.Range("K70").Value = 3902.4
$ws.Range("L70").Value = 2358.75
$ws.Range("M70").Value = -3632.4
$ws.Range("N70").Value = -2898.75
$ws.Range("H73").Value = 1133.919
$ws.Range("I73").Value = 1300.8
$ws.Range("J73").Value = 786.25
$ws.Range("K73").Value = 3902.4
$ws.Range("L73").Value = 2358.75
$ws.Range("M73").Value = -2966.4
$ws.Range("N73").Value = -4230.75
$ws.Range("H105").Value = 24900
$ws.Range("J105").Value = 24900
$ws.Range("L105").Value = 24900
$ws.Range("N105").Value = -31888
$ws.Range("H121").Value = 687.6667
$ws.Range("J121").Value = 626.2
$ws.Range("L121").Value = 1878.6
$ws.Range("N121").Value = -5372.6
$ws.Range("H129").Value = 1805.2941
$ws.Range("J129").Value = 2216.68
$ws.Range("L129").Value = 6650.039999999999
$ws.Range("N129").Value = -16650.04
$ws.Range("H132").Value = 4219.7974
$ws.Range("I132").Value = 3284.0176
$ws.Range("J132").Value = 7357.4116
$ws.Range("K132").Value = 9852.052800000001
$ws.Range("L132").Value = 22072.2348
$ws.Range("M132").Value = -7322.052800000001
$ws.Range("N132").Value = -27132.2348
$ws.Range("H135").Value = 931.63336
$ws.Range("I135").Value = 618.5
$ws.Range("J135").Value = 1792.75
$ws.Range("K135").Value = 5566.5
$ws.Range("L135").Value = 16134.75
$ws.Range("M135").Value = -3031.5
$ws.Range("N135").Value = -21204.75
$ws.Range("H138").Value = 1854.71
$ws.Range("I138").Value = 1345.6531
$ws.Range("J138").Value = 2343.804
$ws.Range("K138").Value = 4036.9593
$ws.Range("L138").Value = 7031.412
$ws.Range("M138").Value = 1103.0407
$ws.Range("N138").Value = -17311.412

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H45").Value = 872.1905
$ws.Range("I45").Value = 866.7778
$ws.Range("J45").Value = 904.6667
$ws.Range("K45").Value = 866.7778
$ws.Range("L45").Value = 904.6667
$ws.Range("M45").Value = -489.7778
$ws.Range("N45").Value = -1658.6667
$ws.Range("H132").Value = 2124380.2
$ws.Range("I132").Value = 5578.909
$ws.Range("J132").Value = 4813628
$ws.Range("K132").Value = 16736.727
$ws.Range("L132").Value = 14440884
$ws.Range("M132").Value = -14206.727
$ws.Range("N132").Value = -14445944

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 838
$ws.Range("I25").Value = 838
$ws.Range("K25").Value = 838
$ws.Range("M25").Value = -603
$ws.Range("H56").Value = 7900
$ws.Range("J56").Value = 7900
$ws.Range("L56").Value = 7900
$ws.Range("N56").Value = -9378
$ws.Range("H109").Value = 29595
$ws.Range("J109").Value = 29595
$ws.Range("L109").Value = 29595
$ws.Range("N109").Value = -32369
$ws.Range("H122").Value = 39120
$ws.Range("J122").Value = 39120
$ws.Range("L122").Value = 39120
$ws.Range("N122").Value = -48920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3333697
$ws.Range("I6").Value = 5000150.5
$ws.Range("J6").Value = 790
$ws.Range("K6").Value = 5000150.5
$ws.Range("L6").Value = 790
$ws.Range("M6").Value = -5000037.5
$ws.Range("N6").Value = -1016
$ws.Range("H22").Value = 640.1905
$ws.Range("I22").Value = 298.22223
$ws.Range("J22").Value = 896.6667
$ws.Range("K22").Value = 298.22223
$ws.Range("L22").Value = 896.6667
$ws.Range("M22").Value = 51.77776999999998
$ws.Range("N22").Value = -1596.6667
$ws.Range("H31").Value = 1650.54
$ws.Range("I31").Value = 1049.6267
$ws.Range("J31").Value = 3453.28
$ws.Range("K31").Value = 1049.6267
$ws.Range("L31").Value = 3453.28
$ws.Range("M31").Value = -754.6267
$ws.Range("N31").Value = -4043.28
$ws.Range("H34").Value = 1650.54
$ws.Range("I34").Value = 1049.6267
$ws.Range("J34").Value = 3453.28
$ws.Range("K34").Value = 1049.6267
$ws.Range("L34").Value = 3453.28
$ws.Range("M34").Value = -847.6267
$ws.Range("N34").Value = -3857.28
$ws.Range("H43").Value = 18916.143
$ws.Range("J43").Value = 18916.143
$ws.Range("L43").Value = 18916.143
$ws.Range("N43").Value = -19284.143
$ws.Range("H50").Value = 22666.666
$ws.Range("J50").Value = 22666.666
$ws.Range("L50").Value = 22666.666
$ws.Range("N50").Value = -23916.666
$ws.Range("H101").Value = 18916.143
$ws.Range("J101").Value = 18916.143
$ws.Range("L101").Value = 18916.143
$ws.Range("N101").Value = -25406.143
$ws.Range("H134").Value = 1907.5217
$ws.Range("I134").Value = 1341.0834
$ws.Range("J134").Value = 2525.4546
$ws.Range("K134").Value = 4023.2502
$ws.Range("L134").Value = 7576.3638
$ws.Range("M134").Value = -1488.2502
$ws.Range("N134").Value = -12646.3638

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 476.5
$ws.Range("I61").Value = 97.59999999999999
$ws.Range("K61").Value = 292.8
$ws.Range("M61").Value = -77.79999999999995
$ws.Range("H113").Value = 827.7778
$ws.Range("I113").Value = 683.3333
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 2049.9999
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = 120.0001000000002
$ws.Range("N113").Value = -7040
$ws.Range("H131").Value = 321248.44
$ws.Range("I131").Value = 445
$ws.Range("J131").Value = 736405.8
$ws.Range("K131").Value = 1335
$ws.Range("L131").Value = 2209217.4
$ws.Range("M131").Value = 3705
$ws.Range("N131").Value = -2219297.4
$ws.Range("H134").Value = 2643.3171
$ws.Range("I134").Value = 1111.6666
$ws.Range("J134").Value = 3842
$ws.Range("K134").Value = 3334.9998
$ws.Range("L134").Value = 11526
$ws.Range("M134").Value = 1735.0002
$ws.Range("N134").Value = -21666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 18052
$ws.Range("J123").Value = 18052
$ws.Range("L123").Value = 18052
$ws.Range("N123").Value = -22952

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1976
$ws.Range("I22").Value = 1950
$ws.Range("J22").Value = 2002
$ws.Range("K22").Value = 1950
$ws.Range("L22").Value = 2002
$ws.Range("M22").Value = -1655
$ws.Range("N22").Value = -2592
$ws.Range("H27").Value = 1976
$ws.Range("I27").Value = 1950
$ws.Range("J27").Value = 2002
$ws.Range("K27").Value = 1950
$ws.Range("L27").Value = 2002
$ws.Range("M27").Value = -1843
$ws.Range("N27").Value = -2216
$ws.Range("H46").Value = 588762.0600000001
$ws.Range("I46").Value = 500.41666
$ws.Range("J46").Value = 2000590
$ws.Range("K46").Value = 500.41666
$ws.Range("L46").Value = 2000590
$ws.Range("M46").Value = -312.41666
$ws.Range("N46").Value = -2000966
$ws.Range("H68").Value = 1861.9048
$ws.Range("J68").Value = 2425
$ws.Range("L68").Value = 2425
$ws.Range("N68").Value = -3923
$ws.Range("H71").Value = 1861.9048
$ws.Range("J71").Value = 2425
$ws.Range("L71").Value = 12125
$ws.Range("N71").Value = -19613
$ws.Range("H111").Value = 31129
$ws.Range("J111").Value = 31129
$ws.Range("L111").Value = 31129
$ws.Range("N111").Value = -39309

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 100011
$ws.Range("J20").Value = 100011
$ws.Range("L20").Value = 100011
$ws.Range("N20").Value = -100491
$ws.Range("H69").Value = 8134.143
$ws.Range("J69").Value = 8134.143
$ws.Range("L69").Value = 8134.143
$ws.Range("N69").Value = -9632.143
$ws.Range("H72").Value = 8134.143
$ws.Range("J72").Value = 8134.143
$ws.Range("L72").Value = 24402.429
$ws.Range("N72").Value = -31890.429
$ws.Range("H113").Value = 857.2857
$ws.Range("I113").Value = 705.25
$ws.Range("J113").Value = 1060
$ws.Range("K113").Value = 2115.75
$ws.Range("L113").Value = 3180
$ws.Range("M113").Value = 54.25
$ws.Range("N113").Value = -7520
$ws.Range("H136").Value = 1272.8064
$ws.Range("I136").Value = 1204.0682
$ws.Range("J136").Value = 1440.8334
$ws.Range("K136").Value = 3612.2046
$ws.Range("L136").Value = 4322.5002
$ws.Range("M136").Value = -1062.2046
$ws.Range("N136").Value = -9422.5002
